$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 416926.75
$ws.Range("I33").Value = 625227.75
$ws.Range("K33").Value = 625227.75
$ws.Range("M33").Value = -624998.75
$ws.Range("H112").Value = 4403.6
$ws.Range("J112").Value = 4403.6
$ws.Range("L112").Value = 13210.8
$ws.Range("N112").Value = -15426.8
$ws.Range("H138").Value = 6782.7676
$ws.Range("J138").Value = 7446.189
$ws.Range("L138").Value = 22338.567
$ws.Range("N138").Value = -32618.567
$ws.Range("H139").Value = 115129.836
$ws.Range("J139").Value = 115129.836
$ws.Range("L139").Value = 115129.836
$ws.Range("N139").Value = -125409.836
$ws.Range("H141").Value = 5989.6924
$ws.Range("I141").Value = 5239.1665
$ws.Range("J141").Value = 14996
$ws.Range("K141").Value = 15717.4995
$ws.Range("L141").Value = 44988
$ws.Range("M141").Value = -10537.4995
$ws.Range("N141").Value = -55348
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H32").Value = 3935.7048
$ws.Range("I32").Value = 2281.5193
$ws.Range("K32").Value = 2281.5193
$ws.Range("M32").Value = -1994.5193
$ws.Range("H45").Value = 3190.3333
$ws.Range("I45").Value = 3190.3333
$ws.Range("K45").Value = 3190.3333
$ws.Range("M45").Value = -2813.3333
$ws.Range("H61").Value = 55224.75
$ws.Range("I61").Value = 51429
$ws.Range("K61").Value = 51429
$ws.Range("M61").Value = -51217
$ws.Range("H110").Value = 4095.6667
$ws.Range("I110").Value = 2325.1428
$ws.Range("J110").Value = 10292.5
$ws.Range("K110").Value = 2325.1428
$ws.Range("L110").Value = 10292.5
$ws.Range("M110").Value = -280.1428000000001
$ws.Range("N110").Value = -14382.5
$ws.Range("H136").Value = 55224.75
$ws.Range("I136").Value = 51429
$ws.Range("K136").Value = 154287
$ws.Range("M136").Value = -151737
$ws.Range("H140").Value = 117000
$ws.Range("J140").Value = 117000
$ws.Range("L140").Value = 117000
$ws.Range("N140").Value = -127360
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 19750
$ws.Range("I36").Value = 19500
$ws.Range("K36").Value = 19500
$ws.Range("M36").Value = -19112
$ws.Range("H40").Value = 19750
$ws.Range("I40").Value = 19500
$ws.Range("K40").Value = 19500
$ws.Range("M40").Value = -19340
$ws.Range("H58").Value = 626928.2
$ws.Range("I58").Value = 1001795.3
$ws.Range("J58").Value = 2149.6667
$ws.Range("K58").Value = 1001795.3
$ws.Range("L58").Value = 2149.6667
$ws.Range("M58").Value = -1001592.3
$ws.Range("N58").Value = -2555.6667
$ws.Range("H99").Value = 9394
$ws.Range("I99").Value = 11637.134
$ws.Range("J99").Value = 6029.3
$ws.Range("K99").Value = 11637.134
$ws.Range("L99").Value = 6029.3
$ws.Range("M99").Value = -10139.134
$ws.Range("N99").Value = -9025.299999999999
$ws.Range("H122").Value = 1923.2609
$ws.Range("I122").Value = 1660.3889
$ws.Range("K122").Value = 4981.1667
$ws.Range("M122").Value = -2531.1667
$ws.Range("H126").Value = 9394
$ws.Range("I126").Value = 11637.134
$ws.Range("J126").Value = 6029.3
$ws.Range("K126").Value = 34911.402
$ws.Range("L126").Value = 18087.9
$ws.Range("M126").Value = -32441.402
$ws.Range("N126").Value = -23027.9
$ws.Range("H132").Value = 8548743
$ws.Range("I132").Value = 10102600
$ws.Range("J132").Value = 2531.3333
$ws.Range("K132").Value = 30307800
$ws.Range("L132").Value = 7593.999899999999
$ws.Range("M132").Value = -30305270
$ws.Range("N132").Value = -12653.9999
$ws.Range("H136").Value = 626928.2
$ws.Range("I136").Value = 1001795.3
$ws.Range("J136").Value = 2149.6667
$ws.Range("K136").Value = 3005385.9
$ws.Range("L136").Value = 6449.000100000001
$ws.Range("M136").Value = -3002835.9
$ws.Range("N136").Value = -11549.0001
$ws.Range("H141").Value = 105685.73
$ws.Range("J141").Value = 120635.555
$ws.Range("L141").Value = 120635.555
$ws.Range("N141").Value = -130995.555
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 55555870
$ws.Range("J23").Value = 62500340
$ws.Range("L23").Value = 187501020
$ws.Range("N23").Value = -187501490
$ws.Range("H29").Value = 125
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 125
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 375
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -929
$ws.Range("H75").Value = 7804.857
$ws.Range("J75").Value = 7804.857
$ws.Range("L75").Value = 23414.571
$ws.Range("N75").Value = -25410.571
$ws.Range("H78").Value = 7804.857
$ws.Range("J78").Value = 7804.857
$ws.Range("L78").Value = 70243.713
$ws.Range("N78").Value = -80227.713
$ws.Range("H80").Value = 5874
$ws.Range("J80").Value = 5874
$ws.Range("L80").Value = 17622
$ws.Range("N80").Value = -19494
$ws.Range("H83").Value = 5874
$ws.Range("J83").Value = 5874
$ws.Range("L83").Value = 52866
$ws.Range("N83").Value = -62226
$ws.Range("H107").Value = 839.3929000000001
$ws.Range("I107").Value = 1527.2858
$ws.Range("J107").Value = 610.0952
$ws.Range("K107").Value = 4581.857400000001
$ws.Range("L107").Value = 1830.2856
$ws.Range("M107").Value = -2661.857400000001
$ws.Range("N107").Value = -5670.2856
$ws.Range("H116").Value = 6001.4614
$ws.Range("I116").Value = 3088.3333
$ws.Range("J116").Value = 8498.429
$ws.Range("K116").Value = 9264.999899999999
$ws.Range("L116").Value = 25495.287
$ws.Range("M116").Value = -5822.999899999999
$ws.Range("N116").Value = -32379.287
$ws.Range("H125").Value = 19648
$ws.Range("I125").Value = 13945
$ws.Range("J125").Value = 22499.5
$ws.Range("K125").Value = 41835
$ws.Range("L125").Value = 67498.5
$ws.Range("M125").Value = -36915
$ws.Range("N125").Value = -77338.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 1500
$ws.Range("I29").Value = 1500
$ws.Range("K29").Value = 1500
$ws.Range("M29").Value = -1210
$ws.Range("H102").Value = 38471690
$ws.Range("I102").Value = 45465812
$ws.Range("K102").Value = 45465812
$ws.Range("M102").Value = -45464190
$ws.Range("H107").Value = 4762642
$ws.Range("J107").Value = 781.1429000000001
$ws.Range("L107").Value = 781.1429000000001
$ws.Range("N107").Value = -4621.1429
$ws.Range("H132").Value = 4414.7393
$ws.Range("I132").Value = 4166.026
$ws.Range("J132").Value = 5800.4287
$ws.Range("K132").Value = 12498.078
$ws.Range("L132").Value = 17401.2861
$ws.Range("M132").Value = -9968.078
$ws.Range("N132").Value = -22461.2861
$ws.Range("H140").Value = 70744.375
$ws.Range("J140").Value = 70744.375
$ws.Range("L140").Value = 70744.375
$ws.Range("N140").Value = -81104.375
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4895.3
$ws.Range("I7").Value = 4120.2085
$ws.Range("K7").Value = 4120.2085
$ws.Range("M7").Value = -4008.2085
$ws.Range("H16").Value = 7694214.5
$ws.Range("I16").Value = 8697090
$ws.Range("K16").Value = 8697090
$ws.Range("M16").Value = -8696920
$ws.Range("H38").Value = 18498.334
$ws.Range("J38").Value = 18498.334
$ws.Range("L38").Value = 18498.334
$ws.Range("N38").Value = -19318.334
$ws.Range("H46").Value = 6454.032
$ws.Range("I46").Value = 4645.1665
$ws.Range("K46").Value = 4645.1665
$ws.Range("M46").Value = -4457.1665
$ws.Range("H126").Value = 4895.3
$ws.Range("I126").Value = 4120.2085
$ws.Range("K126").Value = 12360.6255
$ws.Range("M126").Value = -9890.625499999998
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 199996
$ws.Range("J92").Value = 199996
$ws.Range("L92").Value = 199996
$ws.Range("N92").Value = -204988
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 3500
$ws.Range("J126").Value = 500
$ws.Range("K126").Value = 10500
$ws.Range("L126").Value = 1500
$ws.Range("M126").Value = -8030
$ws.Range("N126").Value = -6440
$ws.Range("H132").Value = 10485754
$ws.Range("I132").Value = 1589037.4
$ws.Range("K132").Value = 4767112.199999999
$ws.Range("M132").Value = -4764582.199999999
